$wb = $excel.ActiveWorkbook

# The workbook has sheets "2024-08-17", "2024-08-18", "2024-08-19" and we
# need to append a new one, "2024-08-20", right after the last one, mirroring
# the existing "archived tasks" sheets (a title in A1, one task per row
# below it).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2024-08-20"

$newSheet.Range("A1").Value = "Archived Tasks"
$newSheet.Range("A2").Value = "Continue 3H Of Python Cours"

# Match the header look of the other sheets' A1 cell (bold text in a
# thin-bordered box) by copying its format from the preceding sheet.
$lastSheet.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# Match the other sheets' page margins (1in top/bottom, 0.5in header/footer)
# instead of the engine's brand-new-sheet defaults.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Adding the new sheet shifts activation/selection onto it; restore the
# workbook's original active sheet/selection so only the sheet addition
# itself is reflected in the saved file.
$originalActive = $wb.Worksheets.Item("2024-08-18")
$originalActive.Activate() | Out-Null
$originalActive.Range("A6").Select() | Out-Null
